# Updates cryptos list figures (prices / 1h volume %) per the latest data pull,
# matching the GitHub Actions scheduled refresh commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Bitcoin)
$ws.Range('D2').Value = '''26.047.47'
$ws.Range('E2').Value = '''  +0.40%  '

# Row 3 (Ethereum)
$ws.Range('D3').Value = '''1.645.58'
$ws.Range('E3').Value = '''  +0.49%  '

# Row 4 (TetherUSD)
$ws.Range('E4').Value = '''  +0.49%  '

# Row 5 (BNB)
$ws.Range('D5').Value = '''215.78'
$ws.Range('E5').Value = '''  +0.50%  '

# Row 6 (XRP)
$ws.Range('E6').Value = '''  +0.14%  '

# Row 7 (USDC)
$ws.Range('E7').Value = '''  +0.50%  '

# Row 8 (Dogecoin)
$ws.Range('E8').Value = '''  +0.70%  '

# Row 9 (Cardano)
$ws.Range('D9').Value = '''0.255'
$ws.Range('E9').Value = '''  +0.47%  '

# Row 10 (Solana)
$ws.Range('D10').Value = '''19.60'
$ws.Range('E10').Value = '''  -0.01%  '

# Row 11 (TRON)
$ws.Range('E11').Value = '''  +0.49%  '

# Row 12 (WrappedEther)
$ws.Range('B12').Value = '''Polkadot'
$ws.Range('C12').Value = '''https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D12').Value = '''4.27'
$ws.Range('E12').Value = '''  +0.56%  '

# Row 13 (Polkadot)
$ws.Range('B13').Value = '''WrappedEther'
$ws.Range('C13').Value = '''https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D13').Value = '''1.613.56'
$ws.Range('E13').Value = '''  -2.08%  '

# Row 14 (Polygon)
$ws.Range('E14').Value = '''  +0.19%  '

# Row 15 (Litecoin)
$ws.Range('D15').Value = '''63.53'
$ws.Range('E15').Value = '''  +1.64%  '

# Row 16 (ShibaInu)
$ws.Range('D16').Value = '''0.0₃0763'
$ws.Range('E16').Value = '''  +0.88%  '

# Row 17 (WrappedBTC)
$ws.Range('D17').Value = '''26.065.52'
$ws.Range('E17').Value = '''  +0.43%  '

# Row 18 (Dai)
$ws.Range('E18').Value = '''  +0.48%  '

# Row 19 (BitcoinCash)
$ws.Range('D19').Value = '''194.54'

# Row 20 (Uniswap)
$ws.Range('D20').Value = '''4.36'
$ws.Range('E20').Value = '''  -0.34%  '

# Row 21 (Avalanche)
$ws.Range('D21').Value = '''9.94'
$ws.Range('E21').Value = '''  +0.21%  '

# Row 22 (Chainlink)
$ws.Range('D22').Value = '''6.21'
$ws.Range('E22').Value = '''  -0.89%  '

# Row 23 (Stellar)
$ws.Range('E23').Value = '''  +4.95%  '

# Row 24 (Toncoin)
$ws.Range('D24').Value = '''1.79'
$ws.Range('E24').Value = '''  +0.01%  '

# Row 25 (BinanceUSD)
$ws.Range('B25').Value = '''Monero'
$ws.Range('C25').Value = '''https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D25').Value = '''144.00'
$ws.Range('E25').Value = '''  -0.12%  '

# Row 26 (Monero)
$ws.Range('B26').Value = '''BinanceUSD'
$ws.Range('C26').Value = '''https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range('D26').Value = '''1.01'
$ws.Range('E26').Value = '''  +0.47%  '

# Row 27 (Cosmos)
$ws.Range('D27').Value = '''6.89'
$ws.Range('E27').Value = '''  +0.67%  '

# Row 28 (EthereumClassic)
$ws.Range('D28').Value = '''15.52'
$ws.Range('E28').Value = '''  +0.38%  '

# Row 29 (PancakeSwap)
$ws.Range('E29').Value = '''  +0.46%  '

# Row 30 (Hedera)
$ws.Range('E30').Value = '''  -1.14%  '

# Row 31 (Filecoin)
$ws.Range('E31').Value = '''  +1.41%  '

# Row 32 (InternetComputer(DFINITY))
$ws.Range('E32').Value = '''  -0.37%  '

# Row 33 (LidoDAOToken)
$ws.Range('E33').Value = '''  -0.14%  '

# Row 34 (HuobiToken)
$ws.Range('D34').Value = '''2.46'
$ws.Range('E34').Value = '''  +1.32%  '

# Row 35 (ARBITRUM)
$ws.Range('D35').Value = '''0.906'
$ws.Range('E35').Value = '''  +0.39%  '

# Row 36 (Maker)
$ws.Range('D36').Value = '''1.131.83'
$ws.Range('E36').Value = '''  -0.50%  '

# Row 37 (ImmutableX)
$ws.Range('E37').Value = '''  -1.05%  '

# Row 39 (VeChain)
$ws.Range('D39').Value = '''0.0157'
$ws.Range('E39').Value = '''  +0.35%  '

# Row 40 (FraxShare)
$ws.Range('E40').Value = '''  +0.83%  '

# Row 41 (Quant)
$ws.Range('D41').Value = '''99.04'
$ws.Range('E41').Value = '''  -0.28%  '

# Row 42 (TrustWalletToken)
$ws.Range('D42').Value = '''0.798'
$ws.Range('E42').Value = '''  -0.12%  '

# Row 43 (BabyDogeCoin)
$ws.Range('E43').Value = '''  +1.89%  '

# Row 44 (Aave)
$ws.Range('D44').Value = '''56.59'

# Row 45 (RenderToken)
$ws.Range('E45').Value = '''  +2.92%  '

# Row 46 (Cronos)
$ws.Range('D46').Value = '''0.0522'
$ws.Range('E46').Value = '''  -1.33%  '

# Row 47 (EnergySwap)
$ws.Range('D47').Value = '''7.80'
$ws.Range('E47').Value = '''  +1.85%  '

# Row 48 (Mantle)
$ws.Range('E48').Value = '''  -0.08%  '

# Row 49 (USDD)
$ws.Range('E49').Value = '''  +0.43%  '

# Row 50 (Algorand)
$ws.Range('E50').Value = '''  -1.05%  '

# Row 51 (NEARProtocol)
$ws.Range('D51').Value = '''1.18'
$ws.Range('E51').Value = '''  +2.88%  '

